# Applies the edits described by the diff to the Mystic Spice Premium
# Chai Tea product description document.

$d = $word.ActiveDocument

function Replace-Exact {
    param(
        [string]$OldText,
        [string]$NewText,
        [bool]$WholeWord = $false
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $true, $WholeWord, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $OldText"
    }
}

# 1) "Principais recursos:" heading becomes bold (w:b val=0 -> w:b, no val)
$rng = $d.Content
$found = $rng.Find.Execute("Principais recursos:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.Bold = 1
} else {
    Write-Host "NOT FOUND: Principais recursos: (for bolding)"
}

# 2) "Mistura" -> "Mistura autêntica"
Replace-Exact "Mistura" "Mistura autêntica" $true

# 3) remainder of the "mistura" bullet
Replace-Exact " autêntica: Nosso chai é uma mistura harmoniosa de folhas de chá preto premium e uma seleção exclusiva de especiarias moídas, incluindo canela, cardamomo, cravo, gengibre e pimenta-do-reino." ": nosso chai é uma mistura harmoniosa de folhas de chá preto premium e uma seleção exclusiva de especiarias moídas, incluindo canela, cardamomo, cravo, gengibre e pimenta preta."

# 4) "Esta receita centenária..." -> Japanese sentence
Replace-Exact "Esta receita centenária promete um sabor autêntico e robusto em cada gole." "この古くから伝わるレシピは、一口飲むごとに本格的でしっかりとした味わいを約束します。"

# 5) "Ingredientes" -> "Ingredientes que melhoram a saúde"
Replace-Exact "Ingredientes" "Ingredientes que melhoram a saúde" $true

# 6) remainder of the "ingredientes" bullet
Replace-Exact " que melhoram a saúde: Cada ingrediente do chá Mystic Spice Chai é escolhido por seus benefícios naturais para a saúde." ": cada ingrediente do Mystic Spice Chai Tea é escolhido por seus benefícios naturais à saúde."

# 7) "Aroma e sabor" -> "Aroma e sabor ricos"
Replace-Exact "Aroma e sabor" "Aroma e sabor ricos" $true

# 8) remainder of the "aroma e sabor" bullet
Replace-Exact " ricos: O aroma quente e picante e o sabor profundo e revigorante do nosso chai fazem dele a bebida perfeita para começar o dia ou relaxar à noite." ": o aroma quente e picante e o sabor profundo e revigorante do nosso chai o tornam a bebida perfeita para começar o dia ou relaxar à noite."

# 9) "Os sabores são intensos..." -> Japanese sentence
Replace-Exact "Os sabores são intensos, mas equilibrados, proporcionando uma experiência reconfortante e relaxante." "風味は強烈でありながらバランスが取れており、快適で心地よい体験を生み出します。"

# 10) "Opções" -> "Opções versáteis de preparo"
Replace-Exact "Opções" "Opções versáteis de preparo" $true

# 11) remainder of the "opções" bullet
Replace-Exact " versáteis de fabricação: Se você ama seu chai fumegante quente, como um chá gelado refrescante ou como um café com leite cremoso, nossa mistura é versátil o suficiente para atender a qualquer preferência." ": não importa se você gosta do seu chai fervendo, como um chá gelado refrescante ou como um café com leite cremoso, nossa mistura é versátil o suficiente para atender a todas as preferências."

# 12) "Instruções simples de preparo..." -> Japanese sentence
Replace-Exact "Instruções simples de preparo estão incluídas para ajudá-lo a saborear seu chai exatamente do jeito que você gosta." "お好みの方法でチャイをお楽しみいただけるよう、簡単な淹れ方の説明書が付属しています。"

# 13) "De" -> "De origem sustentável" (whole word, the standalone heading run)
Replace-Exact "De" "De origem sustentável" $true

# 14) remainder of the "origem sustentável" bullet
Replace-Exact " origem sustentável: Comprometidos com a sustentabilidade, obtemos nossos ingredientes de pequenas fazendas que praticam a agricultura orgânica, garantindo não apenas a melhor qualidade, mas também o bem-estar do nosso planeta." ": comprometidos com a sustentabilidade, adquirimos nossos ingredientes de fazendas pequenas que praticam a agricultura orgânica, garantindo não apenas a melhor qualidade, mas também o bem-estar do nosso planeta."

# 15) "Embalagem" -> "Embalagem elegante"
Replace-Exact "Embalagem" "Embalagem elegante" $true

# 16) remainder of the "embalagem" bullet
Replace-Exact " elegante: O Mystic Spice Chai Tea vem em embalagens ecológicas e com design lindo, tornando-o um presente ideal para os amantes do chá ou um deleite luxuoso para si mesmo." ": O chá Mystic Spice Chai vem em uma embalagem ecológica criada com muita beleza, tornando-a o presente ideal para quem ama chá ou um mimo luxuoso para você."

# 17) "Garantia de Satisfação do Cliente" -> "Garantia de satisfação do cliente"
Replace-Exact "Garantia de Satisfação do Cliente" "Garantia de satisfação do cliente"

# 18) remainder of the "garantia" bullet
Replace-Exact ": Nós apoiamos nosso produto e oferecemos uma garantia de satisfação." ": defendemos nosso produto e oferecemos uma garantia de satisfação."

# 19) "Se o Mystic Spice Chai Tea não atender..." -> Japanese sentence
Replace-Exact "Se o Mystic Spice Chai Tea não atender suas expectativas, estamos comprometidos em resolver da melhor maneira possível." "Mystic Spice Chai Tea がお客様のご期待に添えない場合は、当社が改善するよう努めます。"

# 20) remainder of "ideal para" bullet
Replace-Exact ": Entusiastas do chá, indivíduos preocupados com a saúde, amantes de bebidas quentes e picantes e qualquer pessoa que queira explorar os ricos sabores do tradicional chai indiano." ": entusiastas de chá, pessoas que se preocupam com a saúde, pessoas que adoram bebidas quentes e picantes e qualquer pessoa que queira explorar os ricos sabores do chá indiano tradicional."

Write-Host "Done"
